$d = $word.ActiveDocument
$rng = $d.Content.Duplicate
$rng.Find.Execute("has to", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Found=$($rng.Find.Found) Start=$($rng.Start) End=$($rng.End)"
$rng.Bold = 1
$rng.Text = "must"
$rng.Bold = 0
Write-Output "After edits"
